$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 62.41592166666667
$ws.Range("H2").Value = 187.247765
$ws.Range("I2").Value = 0.1654944774607044
$ws.Range("J2").Value = 0.1654944774607044
$ws.Range("M2").Value = 14.31952633333333
$ws.Range("N2").Value = 42.958579
$ws.Range("O2").Value = 0.7627709815345965
$ws.Range("P2").Value = 0.7627709815345963
$ws.Range("Q2").Value = 893.766433925104
$ws.Range("R2").Value = 8043.897905325935
$ws.Range("S2").Value = 0.1262343850112567
$ws.Range("T2").Value = 0.1262343850112567
$ws.Range("G3").Value = 62.41592166666667
$ws.Range("H3").Value = 187.247765
$ws.Range("I3").Value = 0.1654944774607044
$ws.Range("J3").Value = 0.1654944774607044
$ws.Range("O3").Value = 0.0315532002861196
$ws.Range("P3").Value = 0.0315532002861196
$ws.Range("Q3").Value = 36.97202958863556
$ws.Range("R3").Value = 332.74826629772
$ws.Range("S3").Value = 0.005221880393564314
$ws.Range("T3").Value = 0.005221880393564313
$ws.Range("G4").Value = 62.41592166666667
$ws.Range("H4").Value = 187.247765
$ws.Range("I4").Value = 0.1654944774607044
$ws.Range("J4").Value = 0.1654944774607044
$ws.Range("M4").Value = 1.155932666666667
$ws.Range("N4").Value = 3.467798
$ws.Range("O4").Value = 0.06157409639233436
$ws.Range("P4").Value = 0.06157409639233435
$ws.Range("Q4").Value = 72.14860277460778
$ws.Range("R4").Value = 649.3374249714701
$ws.Range("S4").Value = 0.01019017290756442
$ws.Range("T4").Value = 0.01019017290756442
$ws.Range("G5").Value = 62.41592166666667
$ws.Range("H5").Value = 187.247765
$ws.Range("I5").Value = 0.1654944774607044
$ws.Range("J5").Value = 0.1654944774607044
$ws.Range("M5").Value = 2.705226666666667
$ws.Range("N5").Value = 8.115680000000001
$ws.Range("O5").Value = 0.1441017217869496
$ws.Range("P5").Value = 0.1441017217869495
$ws.Range("Q5").Value = 168.8492157172445
$ws.Range("R5").Value = 1519.6429414552
$ws.Range("S5").Value = 0.02384803914831903
$ws.Range("T5").Value = 0.02384803914831903
$ws.Range("I6").Value = 0.4369365253446571
$ws.Range("J6").Value = 0.436936525344657
$ws.Range("M6").Value = 14.31952633333333
$ws.Range("N6").Value = 42.958579
$ws.Range("O6").Value = 0.7627709815345965
$ws.Range("P6").Value = 0.7627709815345963
$ws.Range("Q6").Value = 2359.711369834961
$ws.Range("R6").Value = 21237.40232851465
$ws.Range("S6").Value = 0.3332825023054602
$ws.Range("T6").Value = 0.3332825023054601
$ws.Range("I7").Value = 0.4369365253446571
$ws.Range("J7").Value = 0.436936525344657
$ws.Range("O7").Value = 0.0315532002861196
$ws.Range("P7").Value = 0.0315532002861196
$ws.Range("S7").Value = 0.01378674569652114
$ws.Range("T7").Value = 0.01378674569652114
$ws.Range("I8").Value = 0.4369365253446571
$ws.Range("J8").Value = 0.436936525344657
$ws.Range("M8").Value = 1.155932666666667
$ws.Range("N8").Value = 3.467798
$ws.Range("O8").Value = 0.06157409639233436
$ws.Range("P8").Value = 0.06157409639233435
$ws.Range("Q8").Value = 190.4858717252016
$ws.Range("R8").Value = 1714.372845526814
$ws.Range("S8").Value = 0.02690397172890356
$ws.Range("T8").Value = 0.02690397172890355
$ws.Range("I9").Value = 0.4369365253446571
$ws.Range("J9").Value = 0.436936525344657
$ws.Range("M9").Value = 2.705226666666667
$ws.Range("N9").Value = 8.115680000000001
$ws.Range("O9").Value = 0.1441017217869496
$ws.Range("P9").Value = 0.1441017217869495
$ws.Range("Q9").Value = 445.793664868249
$ws.Range("R9").Value = 4012.142983814241
$ws.Range("S9").Value = 0.06296330561377222
$ws.Range("T9").Value = 0.0629633056137722
$ws.Range("G10").Value = 57.486235
$ws.Range("H10").Value = 172.458705
$ws.Range("I10").Value = 0.1524235190071549
$ws.Range("J10").Value = 0.1524235190071549
$ws.Range("M10").Value = 14.31952633333333
$ws.Range("N10").Value = 42.958579
$ws.Range("O10").Value = 0.7627709815345965
$ws.Range("P10").Value = 0.7627709815345963
$ws.Range("Q10").Value = 823.1756558866883
$ws.Range("R10").Value = 7408.580902980196
$ws.Range("S10").Value = 0.1162642372020448
$ws.Range("T10").Value = 0.1162642372020447
$ws.Range("G11").Value = 57.486235
$ws.Range("H11").Value = 172.458705
$ws.Range("I11").Value = 0.1524235190071549
$ws.Range("J11").Value = 0.1524235190071549
$ws.Range("O11").Value = 0.0315532002861196
$ws.Range("P11").Value = 0.0315532002861196
$ws.Range("Q11").Value = 34.05193297809333
$ws.Range("R11").Value = 306.46739680284
$ws.Range("S11").Value = 0.004809449823547917
$ws.Range("T11").Value = 0.004809449823547916
$ws.Range("G12").Value = 57.486235
$ws.Range("H12").Value = 172.458705
$ws.Range("I12").Value = 0.1524235190071549
$ws.Range("J12").Value = 0.1524235190071549
$ws.Range("M12").Value = 1.155932666666667
$ws.Range("N12").Value = 3.467798
$ws.Range("O12").Value = 0.06157409639233436
$ws.Range("P12").Value = 0.06157409639233435
$ws.Range("Q12").Value = 66.45021692017667
$ws.Range("R12").Value = 598.05195228159
$ws.Range("S12").Value = 0.009385340451805365
$ws.Range("T12").Value = 0.009385340451805363
$ws.Range("G13").Value = 57.486235
$ws.Range("H13").Value = 172.458705
$ws.Range("I13").Value = 0.1524235190071549
$ws.Range("J13").Value = 0.1524235190071549
$ws.Range("M13").Value = 2.705226666666667
$ws.Range("N13").Value = 8.115680000000001
$ws.Range("O13").Value = 0.1441017217869496
$ws.Range("P13").Value = 0.1441017217869495
$ws.Range("Q13").Value = 155.5132958882667
$ws.Range("R13").Value = 1399.6196629944
$ws.Range("S13").Value = 0.02196449152975686
$ws.Range("T13").Value = 0.02196449152975685
$ws.Range("G14").Value = 92.45614233333333
$ws.Range("H14").Value = 277.368427
$ws.Range("I14").Value = 0.2451454781874835
$ws.Range("J14").Value = 0.2451454781874835
$ws.Range("M14").Value = 14.31952633333333
$ws.Range("N14").Value = 42.958579
$ws.Range("O14").Value = 0.7627709815345965
$ws.Range("P14").Value = 0.7627709815345963
$ws.Range("Q14").Value = 1323.928164820582
$ws.Range("R14").Value = 11915.35348338523
$ws.Range("S14").Value = 0.1869898570158348
$ws.Range("T14").Value = 0.1869898570158348
$ws.Range("G15").Value = 92.45614233333333
$ws.Range("H15").Value = 277.368427
$ws.Range("I15").Value = 0.2451454781874835
$ws.Range("J15").Value = 0.2451454781874835
$ws.Range("O15").Value = 0.0315532002861196
$ws.Range("P15").Value = 0.0315532002861196
$ws.Range("Q15").Value = 54.76633427372177
$ws.Range("R15").Value = 492.897008463496
$ws.Range("S15").Value = 0.007735124372486232
$ws.Range("T15").Value = 0.007735124372486231
$ws.Range("G16").Value = 92.45614233333333
$ws.Range("H16").Value = 277.368427
$ws.Range("I16").Value = 0.2451454781874835
$ws.Range("J16").Value = 0.2451454781874835
$ws.Range("M16").Value = 1.155932666666667
$ws.Range("N16").Value = 3.467798
$ws.Range("O16").Value = 0.06157409639233436
$ws.Range("P16").Value = 0.06157409639233435
$ws.Range("Q16").Value = 106.8730751570829
$ws.Range("R16").Value = 961.857676413746
$ws.Range("S16").Value = 0.01509461130406101
$ws.Range("T16").Value = 0.01509461130406101
$ws.Range("G17").Value = 92.45614233333333
$ws.Range("H17").Value = 277.368427
$ws.Range("I17").Value = 0.2451454781874835
$ws.Range("J17").Value = 0.2451454781874835
$ws.Range("M17").Value = 2.705226666666667
$ws.Range("N17").Value = 8.115680000000001
$ws.Range("O17").Value = 0.1441017217869496
$ws.Range("P17").Value = 0.1441017217869495
$ws.Range("Q17").Value = 250.1148217372622
$ws.Range("R17").Value = 2251.03339563536
$ws.Range("S17").Value = 0.03532588549510147
$ws.Range("T17").Value = 0.03532588549510146
